$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, copying the formatting of the adjacent "sum" header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data values for the "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
